# Adds a new weekly price record for "Ají" (Feria Lagunitas de Puerto Montt)
# by inserting a new row at row 213, which pushes the existing rows
# 213-291 down to 214-292 and extends the used range to A1:R292.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 213 (shifts 213..291 -> 214..292)
$ws.Rows.Item(213).Insert()

# Populate the new row 213 with the new record's data
$ws.Range("A213").Value2 = 4
$ws.Range("B213").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C213").Value2 = "Los Lagos"
$ws.Range("D213").Value2 = 44784
$ws.Range("E213").Value2 = 10
$ws.Range("F213").Value2 = 100112021
$ws.Range("G213").Value2 = "Ají"
$ws.Range("H213").Value2 = "Inferno"
$ws.Range("I213").Value2 = "Primera"
$ws.Range("J213").Value2 = 80
$ws.Range("K213").Value2 = 20000
$ws.Range("L213").Value2 = 20000
$ws.Range("M213").Value2 = 20000
$ws.Range("N213").Value2 = "`$/caja 12 kilos"
$ws.Range("O213").Value2 = "Región de Arica y Parinacota"
$ws.Range("P213").Value2 = 1667
$ws.Range("Q213").Value2 = 12
$ws.Range("R213").Value2 = "Hortaliza"
